$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.042.81'
$ws.Cells.Item(2, 5).Value = '  -0.11%  '
$ws.Cells.Item(3, 4).Value = '1.649.72'
$ws.Cells.Item(3, 5).Value = '  +0.43%  '
$ws.Cells.Item(4, 5).Value = '  -0.31%  '
$ws.Cells.Item(5, 4).Value = '''218.23'
$ws.Cells.Item(5, 5).Value = '  +0.40%  '
$ws.Cells.Item(6, 4).Value = '''0.5191'
$ws.Cells.Item(6, 5).Value = '  +0.30%  '
$ws.Cells.Item(8, 4).Value = '''0.2632'
$ws.Cells.Item(8, 5).Value = '  +0.91%  '
$ws.Cells.Item(9, 4).Value = '''0.06319'
$ws.Cells.Item(9, 5).Value = '  +0.71%  '
$ws.Cells.Item(10, 4).Value = '''20.37'
$ws.Cells.Item(10, 5).Value = '  +0.28%  '
$ws.Cells.Item(11, 4).Value = '''0.07646'
$ws.Cells.Item(12, 4).Value = '''4.584'
$ws.Cells.Item(12, 5).Value = '  +2.66%  '
$ws.Cells.Item(13, 4).Value = '1.644.63'
$ws.Cells.Item(13, 5).Value = '  +0.91%  '
$ws.Cells.Item(14, 4).Value = '1.877.04'
$ws.Cells.Item(15, 4).Value = '''0.5585'
$ws.Cells.Item(15, 5).Value = '  +0.71%  '
$ws.Cells.Item(16, 4).Value = '0.0₅8128'
$ws.Cells.Item(16, 5).Value = '  +1.89%  '
$ws.Cells.Item(17, 4).Value = '''65.18'
$ws.Cells.Item(17, 5).Value = '  +0.85%  '
$ws.Cells.Item(18, 4).Value = '26.025.32'
$ws.Cells.Item(18, 5).Value = '  -0.18%  '
$ws.Cells.Item(19, 5).Value = '  -0.29%  '
$ws.Cells.Item(20, 4).Value = '''4.611'
$ws.Cells.Item(20, 5).Value = '  +0.08%  '
$ws.Cells.Item(21, 5).Value = '  +4.39%  '
$ws.Cells.Item(22, 4).Value = '''191.38'
$ws.Cells.Item(22, 5).Value = '  -0.64%  '
$ws.Cells.Item(23, 5).Value = '  -0.55%  '
$ws.Cells.Item(24, 5).Value = '  -0.32%  '
$ws.Cells.Item(25, 4).Value = '''143.47'
$ws.Cells.Item(25, 5).Value = '  -2.30%  '
$ws.Cells.Item(26, 4).Value = '''0.1184'
$ws.Cells.Item(26, 5).Value = '  -1.52%  '
$ws.Cells.Item(27, 4).Value = '''7.181'
$ws.Cells.Item(27, 5).Value = '  +0.44%  '
$ws.Cells.Item(28, 4).Value = '''15.84'
$ws.Cells.Item(28, 5).Value = '  -0.20%  '
$ws.Cells.Item(29, 4).Value = '''1.508'
$ws.Cells.Item(29, 5).Value = '  +1.99%  '
$ws.Cells.Item(30, 4).Value = '''0.05368'
$ws.Cells.Item(30, 5).Value = '  -4.63%  '
$ws.Cells.Item(31, 5).Value = '  +0.17%  '
$ws.Cells.Item(32, 4).Value = '''3.453'
$ws.Cells.Item(32, 5).Value = '  -0.33%  '
$ws.Cells.Item(33, 4).Value = '''3.345'
$ws.Cells.Item(33, 5).Value = '  -0.43%  '
$ws.Cells.Item(34, 4).Value = '''1.552'
$ws.Cells.Item(34, 5).Value = '  -2.13%  '
$ws.Cells.Item(35, 4).Value = '''2.419'
$ws.Cells.Item(36, 2).Value = 'ARBITRUM'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(36, 4).Value = '''0.9453'
$ws.Cells.Item(36, 5).Value = '  +1.17%  '
$ws.Cells.Item(37, 2).Value = 'MXToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(37, 4).Value = '''2.776'
$ws.Cells.Item(37, 5).Value = '  -0.57%  '
$ws.Cells.Item(38, 4).Value = '''0.5626'
$ws.Cells.Item(38, 5).Value = '  -0.21%  '
$ws.Cells.Item(39, 5).Value = '  +0.23%  '
$ws.Cells.Item(40, 4).Value = '''5.888'
$ws.Cells.Item(40, 5).Value = '  -0.54%  '
$ws.Cells.Item(41, 5).Value = '  -0.23%  '
$ws.Cells.Item(42, 4).Value = '1.029.74'
$ws.Cells.Item(42, 5).Value = '  -2.65%  '
$ws.Cells.Item(43, 4).Value = '''0.8259'
$ws.Cells.Item(43, 5).Value = '  -1.46%  '
$ws.Cells.Item(44, 4).Value = '''100.76'
$ws.Cells.Item(44, 5).Value = '  -1.81%  '
$ws.Cells.Item(45, 4).Value = '1.786.64'
$ws.Cells.Item(45, 5).Value = '  +0.29%  '
$ws.Cells.Item(46, 5).Value = '  +6.35%  '
$ws.Cells.Item(47, 4).Value = '''57.22'
$ws.Cells.Item(47, 5).Value = '  +0.59%  '
$ws.Cells.Item(48, 4).Value = '''1.001'
$ws.Cells.Item(48, 5).Value = '  -0.43%  '
$ws.Cells.Item(49, 4).Value = '''0.4313'
$ws.Cells.Item(49, 5).Value = '  -0.40%  '
$ws.Cells.Item(50, 4).Value = '''7.943'
$ws.Cells.Item(50, 5).Value = '  +0.17%  '
$ws.Cells.Item(51, 5).Value = '  -3.64%  '
